# CRUD-Create Terms Code.xlsx edit
# Commit: Committing System_Setup3 including CRUD testcases of Feature category and VAT Class.
#
# Summary of changes applied:
#  - Sheet "Create_Create Terms Code" (sheet 1): header row C1:K1 relabeled
#    (spaces stripped from the column headings) and selection updated.
#  - Sheet "Edit_Create Terms Code" (sheet 2): PK Terms Code description
#    cells (B2:B4) updated to distinguish rows, column widths widened for
#    columns A, B and D, and selection updated.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Create_Create Terms Code
$ws2 = $wb.Worksheets.Item(2)   # Edit_Create Terms Code

# --- Sheet 2: differentiate the three "PK Terms Code" rows ---
# NOTE: B2 is rewritten first, while its shared string is still referenced
# by only this one cell, so the engine edits that shared-string entry in
# place (keeping its original index) instead of appending a new one.
$ws2.Range("B2").Value = "PK Terms Code Update"

# --- Sheet 1: relabel the header row (spaces removed from each label) ---
$ws1.Range("C1").Value = "DueinDays"
$ws1.Range("D1").Value = "ReceivablesDiscountBaseDateMethod"
$ws1.Range("E1").Value = "DiscountDaysfromBaseDate"
$ws1.Range("F1").Value = "DiscountPercentage"
$ws1.Range("G1").Value = "FinanceChargeBaseDateMethod"
$ws1.Range("H1").Value = "FinanceChargeDaysfromBaseDate"
$ws1.Range("I1").Value = "FinanceChargePercentage"
$ws1.Range("J1").Value = "PayablesDiscountBaseDateMethod`t"
$ws1.Range("K1").Value = "IncTaxandFreightinDiscountableAmt"

# --- Sheet 2: finish differentiating the remaining two rows ---
$ws2.Range("B3").Value = "PK Terms Code Update1"
$ws2.Range("B4").Value = "PK Terms Code Update2"

# --- Sheet 2: widen columns A, B and D to fit the new text ---
$ws2.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws2.Columns.Item(2).ColumnWidth = 19
$ws2.Columns.Item(4).ColumnWidth = 33.5

# --- Update selections on both sheets (sheet 2 stays the active tab) ---
$ws1.Select()
$ws1.Range("A3:K4").Select()

$ws2.Select()
$ws2.Range("B6").Select()
